$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 values (quarterly update of latest period) ---
$ws.Range("B75").Value = 182406
$ws.Range("C75").Value = 6292
$ws.Range("D75").Value = 3003
$ws.Range("E75").Value = 3289
$ws.Range("F75").Value = 25635
$ws.Range("G75").Value = 597
$ws.Range("H75").Value = 25038
$ws.Range("I75").Value = 111
$ws.Range("J75").Value = 6015
$ws.Range("K75").Value = 925
$ws.Range("L75").Value = 5090
$ws.Range("M75").Value = 121149
$ws.Range("N75").Value = 103320
$ws.Range("O75").Value = 1640
$ws.Range("P75").Value = 16190
$ws.Range("Q75").Value = 1496
$ws.Range("R75").Value = 21707
$ws.Range("S75").Value = 4641
$ws.Range("T75").Value = 177765
$ws.Range("U75").Value = 0
$ws.Range("V75").Value = 0
$ws.Range("W75").Value = 10510
$ws.Range("X75").Value = 413
$ws.Range("Y75").Value = 10097
$ws.Range("Z75").Value = 8841
$ws.Range("AA75").Value = 3384
$ws.Range("AB75").Value = 5458
$ws.Range("AC75").Value = 143833
$ws.Range("AD75").Value = 93060
$ws.Range("AE75").Value = 0
$ws.Range("AF75").Value = 50773
$ws.Range("AG75").Value = 14581
$ws.Range("AH75").Value = 39005
$ws.Range("AI75").Value = 45438
$ws.Range("AJ75").Value = 8330
$ws.Range("AK75").Value = 5995

# --- Append new row 76 for period 01-04-2021 ---
# Duplicate row 75 formatting (no explicit styles) down into row 76 first
$ws.Rows.Item(75).Copy()
$ws.Rows.Item(76).PasteSpecial()

# Write the new date label as text (avoid auto date conversion), trimming helper cell afterwards
$ws.Cells.Item(76, 1).Value = "01-04-2021 "
$ws.Cells.Item(76, 2).Formula = "=TRIM(A76)"
$ws.Cells.Item(76, 2).Copy()
$ws.Cells.Item(76, 1).PasteSpecial(-4163)
$ws.Cells.Item(76, 2).Clear()

# Fill in the numeric data for the new row
$ws.Range("B76").Value = 176022
$ws.Range("C76").Value = 7139
$ws.Range("D76").Value = 3532
$ws.Range("E76").Value = 3606
$ws.Range("F76").Value = 23996
$ws.Range("G76").Value = 558
$ws.Range("H76").Value = 23439
$ws.Range("I76").Value = -3
$ws.Range("J76").Value = 5968
$ws.Range("K76").Value = 911
$ws.Range("L76").Value = 5057
$ws.Range("M76").Value = 115376
$ws.Range("N76").Value = 96860
$ws.Range("O76").Value = 1529
$ws.Range("P76").Value = 16986
$ws.Range("Q76").Value = 1392
$ws.Range("R76").Value = 22155
$ws.Range("S76").Value = 2105
$ws.Range("T76").Value = 173917
$ws.Range("U76").Value = 0
$ws.Range("V76").Value = 0
$ws.Range("W76").Value = 10251
$ws.Range("X76").Value = 308
$ws.Range("Y76").Value = 9942
$ws.Range("Z76").Value = 8753
$ws.Range("AA76").Value = 3643
$ws.Range("AB76").Value = 5110
$ws.Range("AC76").Value = 139657
$ws.Range("AD76").Value = 89787
$ws.Range("AE76").Value = 0
$ws.Range("AF76").Value = 49871
$ws.Range("AG76").Value = 15256
$ws.Range("AH76").Value = 37907
$ws.Range("AI76").Value = 41756
$ws.Range("AJ76").Value = 8030
$ws.Range("AK76").Value = 6096
